# Update odds for the first match (row 2) with refreshed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.5
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 2.2
$ws.Range("J2").Value = 4.33
$ws.Range("K2").Value = 1.95
$ws.Range("L2").Value = 3
$ws.Range("Q2").Value = 2.4
$ws.Range("R2").Value = 1.53
$ws.Range("S2").Value = 1.53
$ws.Range("T2").Value = 2.38
$ws.Range("W2").Value = 8
$ws.Range("X2").Value = 15
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 41
$ws.Range("AA2").Value = 34
$ws.Range("AE2").Value = 19
$ws.Range("AG2").Value = 6
$ws.Range("AH2").Value = 9.5
$ws.Range("AI2").Value = 9.5
$ws.Range("AJ2").Value = 21
$ws.Range("AK2").Value = 21
$ws.Range("AN2").Value = 5
$ws.Range("AO2").Value = 21
$ws.Range("AP2").Value = 34
$ws.Range("AQ2").Value = 67
$ws.Range("AT2").Value = 2.38
$ws.Range("AW2").Value = 4
$ws.Range("AX2").Value = 13
$ws.Range("AZ2").Value = 41
$ws.Range("BB2").Value = 251
$ws.Range("BD2").Value = 151

# The match that used to be in row 9 (Progreso vs Fenix, 16:30) no longer
# appears in the weekly schedule; remove its entire row. Excel will shift
# the remaining rows (old row 10 -> new row 9, old row 11 -> new row 10)
# and shrink the used range from BD11 down to BD10 automatically.
$ws.Rows.Item(9).Delete()
